# Edit: "change added to branch1"
#
# Original paragraph 1 (3 runs + proofErr wrapping "git"):
#   "Example " + <proofErr spellStart/> + "git" + <proofErr spellEnd/> + " file."
#   followed by bookmarkStart/bookmarkEnd (_GoBack)
#
# Target:
#   Paragraph 1: single run "Example git file." (no proofErr)
#   Paragraph 2 (new): single run "New change for branch1", followed by the
#                       _GoBack bookmarkStart/bookmarkEnd (moved here)

$d = $word.ActiveDocument

# 1) Re-write the text of paragraph 1 via Find/Replace. Word's Find/Replace
#    normalises the destination range to a single run and clears any
#    proofErr spell-check markers that were splitting the old runs.
$d.Content.Find.Execute("Example git file.", $true, $false, $false, $false, `
    $false, $true, 1, $false, "Example git file.", 2) | Out-Null

# 2) Split paragraph 1 in two, right after its text, creating an empty
#    paragraph 2, then fill it with the new sentence.
$para1Range = $d.Paragraphs(1).Range
$para1Range.Collapse(0)
$para1Range.InsertParagraphAfter()
$para1Range.Collapse(0)
$para1Range.InsertAfter("New change for branch1")

# 3) The new paragraph 2 now holds the text, but the original "_GoBack"
#    bookmark is still sitting at the end of paragraph 1. Temporarily append
#    a one-character sentinel after the new sentence so that "the position
#    right after the text" is an interior run boundary rather than the very
#    edge of the paragraph (collapsed ranges built exactly at a paragraph's
#    trailing edge don't anchor correctly). Use this interior position to
#    move the bookmark, then delete the sentinel again.
$para2End = $d.Paragraphs(2).Range
$para2End = $para2End.Duplicate
$para2End.Collapse(0)
$para2End.InsertAfter("|")

try {
    $d.Bookmarks("_GoBack").Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

$newSentence = $d.Content
$newSentence.Find.Execute("New change for branch1", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkSpot = $newSentence.Duplicate
$bookmarkSpot.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)

$sentinel = $d.Range($bookmarkSpot.Start, $bookmarkSpot.Start + 1)
$sentinel.Delete()
